# Update the predicted LocationX trajectory (column B) with the recomputed
# values coming from the new 3D transformation pipeline that filters out
# z-coordinates below 0. Column A (time) is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    1486.023634431477,
    2298.135986457459,
    3139.354815071962,
    3932.878714179342,
    4565.23899803468,
    5056.847638022701,
    5513.233265149313,
    6063.281889589232,
    6446.602767067534,
    6812.148783421399,
    7149.665591886473,
    7451.723706815204,
    7774.856869573544,
    8127.005400009524,
    8445.84396736661,
    8717.227856540376,
    8889.08285629241,
    9100.810779565025,
    9542.510798792076,
    9699.766527666708,
    9927.335242497451,
    10080.73034756769,
    10379.20822118833,
    10567.4328674492,
    10689.34202809058,
    10874.44348788996,
    10981.21679639662,
    11193.3142242274,
    11323.58291949659,
    11468.94226962042,
    11779.70520856214,
    11864.82703546343,
    12018.08718247714,
    12172.32478077167,
    12245.43284544658,
    12325.09708835992,
    12508.08077212547,
    12635.42560401157,
    12760.11872401384,
    12848.50243680028,
    12981.51286316499,
    13054.52222273246,
    13151.68102437822,
    13287.58764252546,
    13362.2976451646,
    13490.48414440517,
    13593.58868484142,
    13596.94102936735,
    13669.34852991616,
    13776.84786942862,
    13848.21150601052,
    13919.58437113316,
    13963.38835751431,
    14110.50446643876,
    14100.85439618242,
    14189.26885625105,
    14193.50717267355,
    14238.93986774932,
    14387.6725881826,
    14394.71257736065,
    14413.16179364517
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
